$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newLink = "https://api.whatsapp.com/send?phone=6285892254838&text=Halloo%2C%20kak%20NAUFAL%20IMAM%20HIDAYATULLAH%0AAlumni%20UIN%20Syarif%20Hidayatullah%20Jakarta%0AAngkatan%20Wisuda%20131%0A%0AAssalamualaikum%20Wr.Wb.%20%0ADengan%20hormat%2C%20kami%20tim%20Pusat%20Karier%20UIN%20Jakarta%20mengundang%20Anda%20untuk%20berpartisipasi%20dalam%20penyelenggaraan%20Tracer%20Study%20UIN%20Jakarta%202025%20dengan%20mengisi%20survei%20karier%20alumni%20guna%20membantu%20mengevaluasi%20dan%20meningkatkan%20kualitas%20pendidikan%20UIN%20Syarif%20Hidayatullah%20Jakarta.%0A%20%0ABerikut%20merupakan%20langkah-langkah%20untuk%20berpartisipasi%20dalam%20TS%20UIN%20Jakarta%202025.%0ATutorial%20Pengisian%3A%0A1.%20Buka%20Google%20Form%20TS%20Pusat%20Karier%20UIN%20Jakarta%20dengan%20mengakses%20link%20berikut%3A%20https%3A%2F%2Fs.id%2FTSUINJKT2025%0A2.%20Kemudian%2C%20silahkan%20mengisi%20seluruh%20pertanyaan%20dan%20klik%20%E2%80%98Kirim%27%0A3.%20Jika%20ada%20pertanyaan%2C%20silahkan%20berkabar%20melalui%20chat%20WA%20ini%0A%0AAtas%20perhatian%20dan%20partisipasi%20anda%2C%20kami%20ucapkan%20terima%20kasih.%20Semoga%2C%20keberkahan%20dan%20perlindungan%20Allah%20SWT%20selalu%20menyertai%20perjalanan%20karier%20anda.%0A%0ASalam%20Hangat%0A%28Pusat%20Karier%20UIN%20Jakarta%29%0A%0AMengetahui%0AM.%20Kholis%20Hamdy%0A%0AKepala%20Pusat%20Karier%20UIN%20Jakarta"

$ws.Range("A3").Value = "Dummy User 1"
$ws.Range("B3").Value = $newLink
$ws.Range("C3").Value = "Nomor dummy untuk testing parser & bot"

$ws.Range("B2").WrapText = $false

$ws.Range("A5").Select()
